$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("HARSHAD NAGTILAK", 22913.3785, "2025-04"),
    @("HRUTWIK GARDI", 592.861, "2025-04"),
    @("PRATIK RAUL", 29943.4216, "2025-04"),
    @("PRATIK SHIRBHATE", 1752.8626499999998, "2025-04"),
    @("SARANG THAKREY", 557.0264, "2025-04"),
    @("SHUBHAM MUNDADA", 89968.73989999996, "2025-04"),
    @("YUKTA SONIGRA", 16846.567150000003, "2025-04")
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
